$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '48.970.52'
$ws.Range("E2").Value = '  -1.63%  '

# Row 3
$ws.Range("D3").Value = '2.621.72'
$ws.Range("E3").Value = '  +0.41%  '

# Row 4
$ws.Range("E4").Value = '  +0.17%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '110.57'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.55%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '322.47'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.25%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.525'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -1.96%  '

# Row 8
$ws.Range("E8").Value = '  +0.13%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.542'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -3.73%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.53'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -3.59%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.74'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -4.18%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0809'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -1.95%  '

# Row 13
$ws.Range("E13").Value = '  +0.20%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.23'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -1.48%  '

# Row 15
$ws.Range("D15").Value = '3.034.69'
$ws.Range("E15").Value = '  +0.45%  '

# Row 16
$ws.Range("D16").Value = '2.628.32'
$ws.Range("E16").Value = '  +1.39%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.858'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -1.38%  '

# Row 18
$ws.Range("D18").Value = '48.974.05'
$ws.Range("E18").Value = '  -1.55%  '

# Row 19
$ws.Range("E19").Value = '  -4.92%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.80'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -4.08%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.68'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.53%  '

# Row 22
$ws.Range("E22").Value = '  -0.98%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '269.07'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -4.93%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '68.59'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -5.77%  '

# Row 25
$ws.Range("E25").Value = '  -1.04%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.05'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.50%  '

# Row 27
$ws.Range("E27").Value = '  +0.10%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.02'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.58%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.22'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.04%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '34.96'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -3.40%  '

# Row 31
$ws.Range("E31").Value = '  -5.46%  '

# Row 32
$ws.Range("E32").Value = '  -0.50%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.47'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.23%  '

# Row 34
$ws.Range("E34").Value = '  -0.19%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0797'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.33%  '

# Row 36
$ws.Range("B36").Value = 'Celestia'
$ws.Range("C36").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.92'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -4.02%  '

# Row 37
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.99'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +5.30%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.04'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.72%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.12'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.83%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '126.92'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +2.75%  '

# Row 41
$ws.Range("E41").Value = '  -1.66%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '22.28'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -2.67%  '

# Row 43
$ws.Range("E43").Value = '  -4.10%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0317'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.15%  '

# Row 45
$ws.Range("D45").Value = '2.064.54'
$ws.Range("E45").Value = '  +0.87%  '

# Row 46
$ws.Range("B46").Value = 'Stacks'
$ws.Range("C46").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.16'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +6.24%  '

# Row 47
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.24'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -3.05%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.14'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.79%  '

# Row 49
$ws.Range("E49").Value = '  -2.62%  '

# Row 50
$ws.Range("B50").Value = 'MultiversX'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '58.48'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +1.40%  '

# Row 51
$ws.Range("B51").Value = 'THORChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.17'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -3.95%  '
